# eventsliders.xlsx update:
#  - adds labeled alarm sets (new "alarmset(<as>)" command row)
#  - adds Ramp/Soak pattern labels (pidRS can now take a label too)
#  - adds option to load Ramp/Soak patterns from background profile
#  - palette(<int>) can now also be given as a label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# --- pidRS(<int>) -> pidRS(<rs>), with updated description mentioning labels ---
$ws.Range("B76").Value = "pidRS(<rs>)"
$ws.Range("C76").Value = "activates the PID Ramp-Soak pattern number <rs> (1-based!) or the one labeled <rs>"

# --- palette(<int>) -> palette(<p>), with updated description mentioning labels ---
$ws.Range("B83").Value = "palette(<p>)"
$ws.Range("C83").Value = "activates palette <p> with <p> either a number 0-9 or a palette label"

# --- insert a new row above "RC Command" (row 88) for the new alarmset() command ---
$ws.Rows.Item(88).Insert()
$ws.Range("B88").Value = "alarmset(<as>)"
$ws.Range("C88").Value = "activates the alarmset with the given number or label"
